$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
# Row 28
$ws.Range("H28").Value = 6365.353
$ws.Range("I28").Value = 7539
$ws.Range("K28").Value = 7539
$ws.Range("M28").Value = -7054
# Row 64
$ws.Range("H64").Value = 7803.6665
$ws.Range("I64").Value = 3402.875
$ws.Range("J64").Value = 9211.92
$ws.Range("K64").Value = 3402.875
$ws.Range("L64").Value = 9211.92
$ws.Range("M64").Value = -3154.875
$ws.Range("N64").Value = -9707.92
# Row 67
$ws.Range("H67").Value = 7803.6665
$ws.Range("I67").Value = 3402.875
$ws.Range("J67").Value = 9211.92
$ws.Range("K67").Value = 3402.875
$ws.Range("L67").Value = 9211.92
$ws.Range("M67").Value = -2544.875
$ws.Range("N67").Value = -10927.92
# Row 80
$ws.Range("H80").Value = 351
$ws.Range("J80").Value = 465.75
$ws.Range("L80").Value = 1397.25
$ws.Range("N80").Value = -3393.25
# Row 83
$ws.Range("H83").Value = 351
$ws.Range("J83").Value = 465.75
$ws.Range("L83").Value = 4191.75
$ws.Range("N83").Value = -14175.75
# Row 107
$ws.Range("H107").Value = 496
$ws.Range("I107").Value = 312.7857
$ws.Range("K107").Value = 312.7857
$ws.Range("M107").Value = 1607.2143
# Row 125
$ws.Range("H125").Value = 200000
$ws.Range("I125").Value = 0
$ws.Range("J125").Value = 200000
$ws.Range("K125").Value = 0
$ws.Range("L125").Value = 1800000
$ws.Range("M125").ClearContents()
$ws.Range("N125").Value = -1804920
# Row 138
$ws.Range("H138").Value = 2732.087
$ws.Range("J138").Value = 2995.0513
$ws.Range("L138").Value = 8985.153900000001
$ws.Range("N138").Value = -19265.1539

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 3733.2122
$ws.Range("I32").Value = 3662.375
$ws.Range("K32").Value = 3662.375
$ws.Range("M32").Value = -3375.375
# Row 44
$ws.Range("H44").Value = 49990
$ws.Range("J44").Value = 49990
$ws.Range("L44").Value = 49990
$ws.Range("N44").Value = -50966
# Row 55
$ws.Range("H55").Value = 20997
$ws.Range("J55").Value = 49990
$ws.Range("L55").Value = 49990
$ws.Range("N55").Value = -50620
# Row 63
$ws.Range("H63").Value = 5318
$ws.Range("J63").Value = 7366.3335
$ws.Range("L63").Value = 7366.3335
$ws.Range("N63").Value = -8738.333500000001
# Row 66
$ws.Range("H66").Value = 5318
$ws.Range("J66").Value = 7366.3335
$ws.Range("L66").Value = 36831.6675
$ws.Range("N66").Value = -43695.6675
# Row 74
$ws.Range("H74").Value = 2871.5557
$ws.Range("I74").Value = 2871.5557
$ws.Range("K74").Value = 2871.5557
$ws.Range("M74").Value = -1997.5557
# Row 77
$ws.Range("H77").Value = 2871.5557
$ws.Range("I77").Value = 2871.5557
$ws.Range("K77").Value = 14357.7785
$ws.Range("M77").Value = -9989.7785

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
# Row 100
$ws.Range("H100").Value = 73000
$ws.Range("J100").Value = 73000
$ws.Range("L100").Value = 73000
$ws.Range("N100").Value = -75164
# Row 140
$ws.Range("H140").Value = 75000
$ws.Range("J140").Value = 75000
$ws.Range("L140").Value = 75000
$ws.Range("N140").Value = -85360

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Range("H31").Value = 5002.913
$ws.Range("I31").Value = 1598.125
$ws.Range("K31").Value = 1598.125
$ws.Range("M31").Value = -1303.125
# Row 34
$ws.Range("H34").Value = 5002.913
$ws.Range("I34").Value = 1598.125
$ws.Range("K34").Value = 1598.125
$ws.Range("M34").Value = -1396.125
# Row 62
$ws.Range("H62").Value = 8059.143
$ws.Range("I62").Value = 6602.25
$ws.Range("K62").Value = 6602.25
$ws.Range("M62").Value = -5978.25
# Row 65
$ws.Range("H65").Value = 8059.143
$ws.Range("I65").Value = 6602.25
$ws.Range("K65").Value = 33011.25
$ws.Range("M65").Value = -29891.25
# Row 132
$ws.Range("H132").Value = 2316.7334
$ws.Range("I132").Value = 2231.0435
$ws.Range("J132").Value = 2598.2856
$ws.Range("K132").Value = 6693.130500000001
$ws.Range("L132").Value = 7794.8568
$ws.Range("M132").Value = -4163.130500000001
$ws.Range("N132").Value = -12854.8568

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
# Row 13
$ws.Range("H13").Value = 3726.6365
$ws.Range("I13").Value = 332
$ws.Range("K13").Value = 996
$ws.Range("M13").Value = -828
# Row 14
$ws.Range("H14").Value = 18204.723
$ws.Range("I14").Value = 18204.723
$ws.Range("K14").Value = 54614.16900000001
$ws.Range("M14").Value = -54441.16900000001
# Row 26
$ws.Range("H26").Value = 1151.8
$ws.Range("J26").Value = 200
$ws.Range("L26").Value = 600
$ws.Range("N26").Value = -1176

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
# Row 113
$ws.Range("H113").Value = 5673.7393
$ws.Range("I113").Value = 1576.6154
$ws.Range("K113").Value = 1576.6154
$ws.Range("M113").Value = 593.3846000000001
# Row 122
$ws.Range("H122").Value = 5726.2583
$ws.Range("I122").Value = 5778.5557
$ws.Range("J122").Value = 5373.25
$ws.Range("K122").Value = 17335.6671
$ws.Range("L122").Value = 16119.75
$ws.Range("M122").Value = -14885.6671
$ws.Range("N122").Value = -21019.75

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
# Row 40
$ws.Range("H40").Value = 6160.4
$ws.Range("I40").Value = 4459.4165
$ws.Range("J40").Value = 7730.5386
$ws.Range("K40").Value = 4459.4165
$ws.Range("L40").Value = 7730.5386
$ws.Range("M40").Value = -4323.4165
$ws.Range("N40").Value = -8002.5386
# Row 46
$ws.Range("H46").Value = 2830.4614
$ws.Range("I46").Value = 800
$ws.Range("J46").Value = 2999.6667
$ws.Range("K46").Value = 800
$ws.Range("L46").Value = 2999.6667
$ws.Range("M46").Value = -612
$ws.Range("N46").Value = -3375.6667
# Row 63
$ws.Range("H63").Value = 25849.857
$ws.Range("I63").Value = 20316.666
$ws.Range("J63").Value = 29999.75
$ws.Range("K63").Value = 20316.666
$ws.Range("L63").Value = 29999.75
$ws.Range("M63").Value = -19567.666
$ws.Range("N63").Value = -31497.75
# Row 66
$ws.Range("H66").Value = 25849.857
$ws.Range("I66").Value = 20316.666
$ws.Range("J66").Value = 29999.75
$ws.Range("K66").Value = 60949.99800000001
$ws.Range("L66").Value = 89999.25
$ws.Range("M66").Value = -57205.99800000001
$ws.Range("N66").Value = -97487.25
# Row 122
$ws.Range("H122").Value = 6164.222
$ws.Range("I122").Value = 4746.6665
$ws.Range("J122").Value = 8999.333000000001
$ws.Range("K122").Value = 14239.9995
$ws.Range("L122").Value = 26997.999
$ws.Range("M122").Value = -11789.9995
$ws.Range("N122").Value = -31897.999

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
# Row 70
$ws.Range("H70").Value = 27500
$ws.Range("I70").Value = 23333.334
$ws.Range("J70").Value = 40000
$ws.Range("K70").Value = 23333.334
$ws.Range("L70").Value = 40000
$ws.Range("M70").Value = -23018.334
$ws.Range("N70").Value = -40630
# Row 73
$ws.Range("H73").Value = 27500
$ws.Range("I73").Value = 23333.334
$ws.Range("J73").Value = 40000
$ws.Range("K73").Value = 23333.334
$ws.Range("L73").Value = 40000
$ws.Range("M73").Value = -22241.334
$ws.Range("N73").Value = -42184
# Row 122
$ws.Range("H122").Value = 6291.5
$ws.Range("I122").Value = 1526.4286
$ws.Range("J122").Value = 9997.666999999999
$ws.Range("K122").Value = 4579.2858
$ws.Range("L122").Value = 29993.001
$ws.Range("M122").Value = -2129.2858
$ws.Range("N122").Value = -34893.001
# Row 136
$ws.Range("H136").Value = 2078.125
$ws.Range("J136").Value = 3409
$ws.Range("L136").Value = 10227
$ws.Range("N136").Value = -15327
